# Refresh Sheet3's lookup table (B20:B36) with the latest coefficients.
# Sheet1!CB:CC hold VLOOKUP()s against this table, so they recalc in place.
$wb = $excel.ActiveWorkbook
$sheet3 = $wb.Worksheets.Item("Sheet3")

$sheet3.Range("B20").Value = 14.11258681578026
$sheet3.Range("B21").Value = 2.3061734304100252
$sheet3.Range("B22").Value = 8.6477451095384463
$sheet3.Range("B23").Value = 4.9420866472597345
$sheet3.Range("B24").Value = 9.5146973188374702
$sheet3.Range("B26").Value = 4.5359485055591584
$sheet3.Range("B27").Value = 7.7016289315344393
$sheet3.Range("B28").Value = 6.1097237465596645
$sheet3.Range("B29").Value = 5.6087545349707479
$sheet3.Range("B30").Value = 0.071876271447111131
$sheet3.Range("B31").Value = 1.0375768621809889
$sheet3.Range("B32").Value = 6.6542717890275496
$sheet3.Range("B33").Value = 4.4376056965310759
$sheet3.Range("B34").Value = 0.53256600935462006
$sheet3.Range("B35").Value = 9.3189498171769962
$sheet3.Range("B36").Value = 43.826786251136497

$null = $excel.Calculate()

# Sheet1 gains a new "as of" snapshot column (CF, "01-nov") - the day after
# the existing last snapshot column (CE, "31-oct") - holding the freshly
# recalculated VLOOKUP results as static values, the same way CE was built
# from CB/CC.
$sheet1 = $wb.Worksheets.Item("Sheet1")

$sheet1.Range("CF1").Value = "01-nov"
$sheet1.Range("CF1").NumberFormat = $sheet1.Range("CE1").NumberFormat

for ($r = 2; $r -le 18; $r++) {
    $cb = $sheet1.Cells.Item($r, 80).Value()
    $sheet1.Cells.Item($r, 84).Value = $cb
    $sheet1.Cells.Item($r, 84).NumberFormat = $sheet1.Cells.Item($r, 83).NumberFormat
}

$null = $sheet1.Range("CF3").Select()
